$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.374.02'
$ws.Range('E2').Value = '  +3.00%  '
$ws.Range('D3').Value = '2.413.44'
$ws.Range('E3').Value = '  +0.78%  '
$ws.Range('E4').Value = '  +0.30%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '572.88'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.36%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.61'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.55%  '
$ws.Range('E7').Value = '  -0.32%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.538'
$ws.Range('D8').NumberFormat = 'General'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.29%  '
$ws.Range('D9').Value = '2.437.60'
$ws.Range('E9').Value = '  +1.61%  '
$ws.Range('E10').Value = '  +4.64%  '
$ws.Range('E11').Value = '  +0.66%  '
$ws.Range('E12').Value = '  +3.88%  '
$ws.Range('E13').Value = '  +3.59%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.87'
$ws.Range('D14').NumberFormat = 'General'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +4.40%  '
$ws.Range('E15').Value = '  +8.06%  '
$ws.Range('D16').Value = '2.883.50'
$ws.Range('E16').Value = '  +1.28%  '
$ws.Range('D17').Value = '62.130.22'
$ws.Range('E17').Value = '  +2.63%  '
$ws.Range('D18').Value = '2.441.17'
$ws.Range('E18').Value = '  +1.97%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.91'
$ws.Range('D19').NumberFormat = 'General'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -5.24%  '
$ws.Range('E20').Value = '  +2.49%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '326.12'
$ws.Range('D21').NumberFormat = 'General'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.27%  '
$ws.Range('E22').Value = '  +2.23%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.04'
$ws.Range('D23').NumberFormat = 'General'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +13.74%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.00'
$ws.Range('D24').NumberFormat = 'General'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.09%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '65.44'
$ws.Range('D25').NumberFormat = 'General'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.46%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '615.84'
$ws.Range('D26').NumberFormat = 'General'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +10.29%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.40'
$ws.Range('D27').NumberFormat = 'General'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +4.28%  '
$ws.Range('D28').Value = '0.0₃0986'
$ws.Range('E28').Value = '  +8.06%  '
$ws.Range('D29').Value = '2.531.78'
$ws.Range('E29').Value = '  +0.37%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.997'
$ws.Range('D30').NumberFormat = 'General'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.06%  '
$ws.Range('E31').Value = '  +1.99%  '
$ws.Range('E32').Value = '  +8.19%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.137'
$ws.Range('D33').NumberFormat = 'General'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +5.69%  '
$ws.Range('E34').Value = '  +1.80%  '
$ws.Range('E35').Value = '  +4.50%  '
$ws.Range('E36').Value = '  -0.37%  '
$ws.Range('E37').Value = '  +5.16%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '152.97'
$ws.Range('D38').NumberFormat = 'General'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.36%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.373'
$ws.Range('D39').NumberFormat = 'General'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.20%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.40'
$ws.Range('D40').NumberFormat = 'General'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +6.51%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '18.56'
$ws.Range('D41').NumberFormat = 'General'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.72%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.73'
$ws.Range('D42').NumberFormat = 'General'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +18.02%  '
$ws.Range('E43').Value = '  +6.77%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '42.31'
$ws.Range('D44').NumberFormat = 'General'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.89%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.999'
$ws.Range('D45').NumberFormat = 'General'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.01%  '
$ws.Range('D46').Value = '0.0₆0283'
$ws.Range('E46').Value = '  +0.66%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '143.78'
$ws.Range('D47').NumberFormat = 'General'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.43%  '
$ws.Range('E48').Value = '  +2.53%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '20.25'
$ws.Range('D49').NumberFormat = 'General'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +6.59%  '
$ws.Range('E50').Value = '  +2.21%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0514'
$ws.Range('D51').NumberFormat = 'General'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.58%  '
